# Adds 4 new incident rows (145-148) to the bottom of the "v2_management"
# log sheet, matching the "updating errors / fixing issues / v.5.5" commit.
#
# All cells in this sheet are plain text (Bloque / Incidencia / Fecha / Hora /
# Turno / Hora de Reparacion / Tiempo de Reparacion / MTBF columns), so the
# "Fecha" column (C) needs to be pre-formatted as Text before the value is
# written - otherwise Excel's normal type-inference on Range.Value would
# silently turn a literal "2024-06-10" string into a real date serial
# number, which would not match the source data (a plain string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(145, "WC48 P5F", "AOI (malla)",                     "2024-06-10", "09:46:28", "Mañana", "09:46:29", "0:00:01", "N/A"),
    @(146, "WC48 P5F", "Cámara no detecta foams",          "2024-06-10", "09:46:41", "Mañana", "09:46:42", "0:00:01", "N/A"),
    @(147, "WC48 P5F", "Cámara no detecta foam derecho",   "2024-06-10", "09:46:44", "Mañana", "09:46:45", "0:00:01", "N/A"),
    @(148, "WC48 P5F", "Cámara no detecta skeleton",       "2024-06-10", "09:46:48", "Mañana", "09:46:50", "0:00:02", "N/A")
)

foreach ($row in $rows) {
    $r = $row[0]

    # Force column C ("Fecha") to Text so the date-like string is stored
    # verbatim instead of being coerced into a date serial number, then put
    # the cell style back to "Normal" so it ends up with no special
    # formatting left behind (matching every other cell in the sheet).
    $ws.Cells.Item($r, 3).NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]

    $ws.Cells.Item($r, 3).Style = "Normal"
}
